$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3: fill in new hac/fast columns (L3, K3, J3) - order matches shared-string append order
$ws.Range("L3").Value = "/projects/MicroBench/data/PAS01578.dorado0.7.3.bmdna_r10.4.1_e8.2_400bps@5.0.0_hac.dup.fastq.gz"
$ws.Range("K3").Value = "/projects/MicroBench/data/PAS01578.dorado0.7.3.bmdna_r10.4.1_e8.2_400bps@5.0.0_hac.sim.fastq.gz"
$ws.Range("J3").Value = """/projects/MicroBench/data/.PAS01578.dorado0.7.3.bmdna_r10.4.1_e8.2_400bps@5.0.0_fast.sim.fastq.gz"""

# Header O1: " allmods" -> " mods"
$ws.Range("O1").Value = " mods"

# Update selected cell to O2 as shown in the diff
$ws.Range("O2").Select()
